# The commit rotates the species-record rows 2-5 (columns A,B,D,E,F,G,H,Q,R)
# cyclically: row2 gets what row3 had, row3 gets what row4 had, row4 gets
# what row5 had, and row5 gets what row2 originally had. Column C (always
# "Ovaliderad") is unaffected and thus untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")
$rows = @(2, 3, 4, 5)

# Snapshot the original values for the affected columns/rows before writing
# anything, since the rotation reads from rows we are about to overwrite.
$orig = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $orig[$addr] = $ws.Range($addr).Value2
    }
}

# New row order: row N takes the old values of row N+1 (wrapping 5 -> 2).
$srcRowFor = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($r in $rows) {
    $srcRow = $srcRowFor[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $orig["$col$srcRow"]
    }
}
